$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range("D2").Value = "34.171.95"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.799.65"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "2.059.57"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("E14").Value = "  +10.74%  "
$ws.Range("D15").Value = "1.803.08"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "34.189.76"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  +3.36%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("D36").Value = "1.504.01"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "1.956.79"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("E49").Value = "  +7.42%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  -5.22%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E45").Value = "  -2.18%  "

# --- Values that look numeric (e.g. '1.00', '0.559'): force text so
#     Excel does not coerce them to a Double and lose formatting ---
$numberLikeRefs = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D14","D16","D18","D19","D20","D23","D26","D27","D30","D31","D38","D39","D41","D43","D48","D49","D50","D51","D44","D45")
foreach ($r in $numberLikeRefs) {
    $ws.Range($r).NumberFormat = "@"
}
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "227.77"
$ws.Range("D6").Value = "0.559"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").Value = "31.52"
$ws.Range("D9").Value = "46.19"
$ws.Range("D10").Value = "0.283"
$ws.Range("D11").Value = "0.0665"
$ws.Range("D12").Value = "0.0928"
$ws.Range("D14").Value = "11.49"
$ws.Range("D16").Value = "0.641"
$ws.Range("D18").Value = "4.25"
$ws.Range("D19").Value = "70.13"
$ws.Range("D20").Value = "255.28"
$ws.Range("D23").Value = "10.44"
$ws.Range("D26").Value = "158.77"
$ws.Range("D27").Value = "16.67"
$ws.Range("D30").Value = "1.00"
$ws.Range("D31").Value = "3.96"
$ws.Range("D38").Value = "0.639"
$ws.Range("D39").Value = "86.20"
$ws.Range("D41").Value = "2.84"
$ws.Range("D43").Value = "0.912"
$ws.Range("D48").Value = "5.76"
$ws.Range("D49").Value = "11.99"
$ws.Range("D50").Value = "1.00"
$ws.Range("D51").Value = "51.79"
$ws.Range("D44").Value = "0.0520"
$ws.Range("D45").Value = "2.06"
foreach ($r in $numberLikeRefs) {
    $ws.Range($r).Style = "Normal"
}
